$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update the Date value (row 8, column B: "Date" label is in A8) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Sheet "Include #0": update System URI for TRE-R43-CapaciteSavoirFaire ---
$inc0 = $wb.Worksheets.Item("Include #0")
$inc0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R43-CapaciteSavoirFaire"

# --- Sheet "Include #1": update System URI for TRE-R01-EnsembleSavoirFaire-CISIS ---
$inc1 = $wb.Worksheets.Item("Include #1")
$inc1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R01-EnsembleSavoirFaire-CISIS"
